$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

$newRow = 15

$ws.Cells.Item($newRow, 1).Value = "Kun jij dit even regelen?"
$ws.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$ws.Cells.Item($newRow, 3).Value = "Testmail #1: Kun jij dit even regelen?"
$ws.Cells.Item($newRow, 4).Value = "Planning / Afspraak"
$ws.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$ws.Cells.Item($newRow, 6).Value = "2025-08-05 17:15:24"
$ws.Cells.Item($newRow, 7).Value = "Ja"
$ws.Cells.Item($newRow, 8).Value = "Ja"
$ws.Cells.Item($newRow, 9).Value = "Nee"
$ws.Cells.Item($newRow, 10).Value = "Nee"

# Update conditional formatting ranges to include the new row
$ws.Range("D2:D15").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D15"))
$ws.Range("G2:G15").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("G2:G15"))
$ws.Range("H2:H15").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("H2:H15"))
$ws.Range("I2:I15").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("I2:I15"))
$ws.Range("J2:J15").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("J2:J15"))

# Update the Dashboard summary count for "Planning / Afspraak"
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 9
